$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the four date cells back by one day (feedback processing update)
$ws.Range("A7").Value = 42184
$ws.Range("A10").Value = 42185
$ws.Range("A13").Value = 42186
$ws.Range("A16").Value = 42187
